# Convert the "Order Amount" column (D) of the sales table from text
# dollar-formatted strings (e.g. "$925.00") to real numbers formatted
# with a custom Rupee currency number format ("₹"#,##0.00).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New raw numeric values for D2:D40 (Order Amount), row-by-row.
$orderAmounts = @(
  2174, 10469, 30283, 8536,
  19757, 12890, 37542, 29218, 4811,
  20937, 33659, 7594, 18361,
  29585, 10863, 24492,
  6257, 31905, 13927, 27189, 1624, 38777,
  9533, 21815, 12641, 35097,
  2920, 17364, 4087, 29681, 6756,
  21349, 14011, 27005, 9572, 30259, 16747, 3848, 25060
)

$firstRow = 2
for ($i = 0; $i -lt $orderAmounts.Length; $i++) {
  $row = $firstRow + $i
  $ws.Cells.Item($row, 4).Value = $orderAmounts[$i]
}

# Apply the custom currency number format to the whole data column.
$lo = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Item(4)
$col.DataBodyRange.NumberFormat = """₹""#,##0.00"

Write-Output "Order Amount column converted to numeric currency values"
